# Rename the diff-comparison column headers so the "old"/"new" suffixes
# reflect the actual format-version names (FV2404 / FV2410), then turn the
# data range into a real Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "<name>_old" -> "<name>_FV2404"
for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseHeaders[$i] + "_FV2404"
}

# Column K (11) stays "diff" - untouched.

# Columns L-U (12-21): "<name>_new" -> "<name>_FV2410"
for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseHeaders[$i] + "_FV2410"
}

# Turn A1:U64 into a proper Excel Table ("Table1") with a header row.
$dataRange = $ws.Range("A1:U64")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (split below row 1).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
